$wb = $excel.ActiveWorkbook

# ALC row 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 88.40000000000001
$ws.Range("I4").Value = 108
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 108
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = -238

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 645
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 754
$ws.Range("K8").Value = 300
$ws.Range("L8").Value = 2262
$ws.Range("M8").Value = -161
$ws.Range("N8").Value = -2540

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 83334984
$ws.Range("I100").Value = 1850
$ws.Range("J100").Value = 250001250
$ws.Range("K100").Value = 1850
$ws.Range("L100").Value = 250001250
$ws.Range("M100").Value = -1309
$ws.Range("N100").Value = -250002332

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2429.127
$ws.Range("I138").Value = 2342.625
$ws.Range("J138").Value = 2482.359
$ws.Range("K138").Value = 7027.875
$ws.Range("L138").Value = 7447.076999999999
$ws.Range("M138").Value = -1887.875
$ws.Range("N138").Value = -17727.077

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2282.7368
$ws.Range("I141").Value = 2240.8572
$ws.Range("J141").Value = 2400
$ws.Range("K141").Value = 6722.571599999999
$ws.Range("L141").Value = 7200
$ws.Range("M141").Value = -1542.571599999999
$ws.Range("N141").Value = -17560

# ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 37900
$ws.Range("J80").Value = 37900
$ws.Range("L80").Value = 37900
$ws.Range("N80").Value = -39896

# ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 37900
$ws.Range("J83").Value = 37900
$ws.Range("L83").Value = 113700
$ws.Range("N83").Value = -123684

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 69040.766
$ws.Range("J82").Value = 36415.715
$ws.Range("L82").Value = 36415.715
$ws.Range("N82").Value = -37181.715

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 69040.766
$ws.Range("J85").Value = 36415.715
$ws.Range("L85").Value = 36415.715
$ws.Range("N85").Value = -39067.715

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1760.3334
$ws.Range("I99").Value = 1407.5
$ws.Range("J99").Value = 1821.6957
$ws.Range("K99").Value = 1407.5
$ws.Range("L99").Value = 1821.6957
$ws.Range("M99").Value = 90.5
$ws.Range("N99").Value = -4817.6957

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 18009.166
$ws.Range("I41").Value = 4900
$ws.Range("J41").Value = 20631
$ws.Range("K41").Value = 4900
$ws.Range("L41").Value = 20631
$ws.Range("M41").Value = -4472
$ws.Range("N41").Value = -21487

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9363.333000000001
$ws.Range("J51").Value = 9363.333000000001
$ws.Range("L51").Value = 9363.333000000001
$ws.Range("N51").Value = -10835.333

# CRP row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 26061
$ws.Range("J60").Value = 26061
$ws.Range("L60").Value = 26061
$ws.Range("N60").Value = -27083

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9363.333000000001
$ws.Range("J61").Value = 9363.333000000001
$ws.Range("L61").Value = 9363.333000000001
$ws.Range("N61").Value = -10059.333

# CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 26303.95
$ws.Range("J88").Value = 26907.578
$ws.Range("L88").Value = 26907.578
$ws.Range("N88").Value = -27719.578

# CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 26303.95
$ws.Range("J91").Value = 26907.578
$ws.Range("L91").Value = 26907.578
$ws.Range("N91").Value = -29715.578

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1966
$ws.Range("I99").Value = 1783.3334
$ws.Range("J99").Value = 2240
$ws.Range("K99").Value = 1783.3334
$ws.Range("L99").Value = 2240
$ws.Range("M99").Value = -285.3334
$ws.Range("N99").Value = -5236

# CRP row 109
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 11900
$ws.Range("J109").Value = 11900
$ws.Range("L109").Value = 11900
$ws.Range("N109").Value = -13980

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1966
$ws.Range("I126").Value = 1783.3334
$ws.Range("J126").Value = 2240
$ws.Range("K126").Value = 5350.0002
$ws.Range("L126").Value = 6720
$ws.Range("M126").Value = -2880.0002
$ws.Range("N126").Value = -11660

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 12359
$ws.Range("J57").Value = 16011
$ws.Range("L57").Value = 16011
$ws.Range("N57").Value = -17651

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 44457.883
$ws.Range("I70").Value = 49586.824
$ws.Range("J70").Value = 5136
$ws.Range("K70").Value = 49586.824
$ws.Range("L70").Value = 5136
$ws.Range("M70").Value = -49316.824
$ws.Range("N70").Value = -5676

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 44457.883
$ws.Range("I73").Value = 49586.824
$ws.Range("J73").Value = 5136
$ws.Range("K73").Value = 49586.824
$ws.Range("L73").Value = 5136
$ws.Range("M73").Value = -48650.824
$ws.Range("N73").Value = -7008

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3899.2307
$ws.Range("I97").Value = 3899.2307
$ws.Range("K97").Value = 3899.2307
$ws.Range("M97").Value = -3403.2307

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1524
$ws.Range("I102").Value = 940.4
$ws.Range("J102").Value = 2357.7144
$ws.Range("K102").Value = 940.4
$ws.Range("L102").Value = 2357.7144
$ws.Range("M102").Value = 681.6
$ws.Range("N102").Value = -5601.7144

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1907.4445
$ws.Range("I93").Value = 1728.6666
$ws.Range("J93").Value = 2801.3333
$ws.Range("K93").Value = 1728.6666
$ws.Range("L93").Value = 2801.3333
$ws.Range("M93").Value = -480.6666
$ws.Range("N93").Value = -5297.3333

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1799.5
$ws.Range("I100").Value = 1799.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1799.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1258.5
$ws.Range("N100").ClearContents()

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2981.52
$ws.Range("I136").Value = 2279.889
$ws.Range("J136").Value = 4785.7144
$ws.Range("K136").Value = 6839.667
$ws.Range("L136").Value = 14357.1432
$ws.Range("M136").Value = -4289.667
$ws.Range("N136").Value = -19457.1432

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2889.5789
$ws.Range("I107").Value = 1709.3636
$ws.Range("J107").Value = 4512.375
$ws.Range("K107").Value = 5128.0908
$ws.Range("L107").Value = 13537.125
$ws.Range("M107").Value = -3208.0908
$ws.Range("N107").Value = -17377.125

# WVR row 109
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 26318
$ws.Range("J109").Value = 26318
$ws.Range("L109").Value = 26318
$ws.Range("N109").Value = -29092

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 126777.375
$ws.Range("I113").Value = 143174.14
$ws.Range("K113").Value = 429522.42
$ws.Range("M113").Value = -427352.42

Write-Host "Applied all changes successfully"
